$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "245.16"
Set-TextValue "D3" "29.08"
Set-TextValue "E3" "-1.92%"
Set-TextValue "D4" "5.255"
Set-TextValue "E4" "1.52%"
Set-TextValue "D5" "0.05703"
Set-TextValue "E5" "-0.10%"
Set-TextValue "E6" "0.17%"
Set-TextValue "D7" "3.193"
Set-TextValue "E7" "3.88%"
Set-TextValue "D8" "0.8517"
Set-TextValue "E8" "-0.57%"
Set-TextValue "D9" "0.8544"
Set-TextValue "E9" "-1.85%"
Set-TextValue "E10" "0.27%"
Set-TextValue "D11" "0.07049"
Set-TextValue "E11" "-0.53%"
Set-TextValue "D12" "0.03192"
Set-TextValue "E12" "9.22%"
Set-TextValue "D13" "0.09248"
Set-TextValue "E13" "-1.47%"
Set-TextValue "D14" "0.001527"
Set-TextValue "E14" "0.99%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005996"
Set-TextValue "E15" "-2.36%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.493"
Set-TextValue "E16" "0.20%"
Set-TextValue "B17" "BTSEToken"
Set-TextValue "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D17" "2.175"
Set-TextValue "E17" "-0.55%"
Set-TextValue "B18" "One"
Set-TextValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005954"
Set-TextValue "E18" "-1.11%"
Set-TextValue "D19" "0.3158"
Set-TextValue "E19" "-0.50%"
Set-TextValue "D20" "0.03270"
Set-TextValue "E20" "-2.27%"
Set-TextValue "E21" "-1.95%"
Set-TextValue "D22" "3.491"
Set-TextValue "E22" "0.64%"
Set-TextValue "D23" "0.04089"
Set-TextValue "E23" "-2.13%"
Set-TextValue "E24" "0.01%"
Set-TextValue "E25" "0.05%"
Set-TextValue "D26" "0.004139"
Set-TextValue "E26" "-17.65%"
Set-TextValue "D27" "0.0001201"
Set-TextValue "D28" "0.0001450"
Set-TextValue "E28" "-25.21%"
Set-TextValue "D40" "0.03754"
Set-TextValue "E40" "0.23%"
Set-TextValue "D41" "0.1063"
Set-TextValue "E41" "-0.74%"
Set-TextValue "D42" "0.003709"
Set-TextValue "E42" "6.64%"
Set-TextValue "D43" "0.002403"
Set-TextValue "E43" "-5.41%"
Set-TextValue "D44" "0.009365"
Set-TextValue "E44" "2.03%"
Set-TextValue "D45" "0.00005295"
Set-TextValue "E45" "1.57%"
Set-TextValue "E46" "0.05%"
Set-TextValue "D47" "0.07505"
Set-TextValue "E47" "29.38%"
Set-TextValue "D49" "0.00002101"
Set-TextValue "E49" "0.05%"
Set-TextValue "D50" "0.0002001"
Set-TextValue "E50" "0.05%"

Write-Host "Applied $($ws.UsedRange.Rows.Count) row updates"
